$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear all existing content/formatting ---
$ws.Cells.Clear() | Out-Null

# --- Column widths (best-effort match to target widths given engine pixel quantization) ---
$ws.Columns.Item(2).ColumnWidth = 8.0
$ws.Columns.Item(3).ColumnWidth = 35.166666666666664
$ws.Columns.Item(4).ColumnWidth = 51.0
$ws.Columns.Item(5).ColumnWidth = 34.166666666666664
$ws.Columns.Item(6).ColumnWidth = 28.5
$ws.Columns.Item(7).ColumnWidth = 9.333333333333334
$ws.Columns.Item(8).ColumnWidth = 7.833333333333333

# --- Cell values ---
$ws.Range("A1").Value = "TestID"
$ws.Range("B1").Value = "Enabled"
$ws.Range("C1").Value = "action"
$ws.Range("D1").Value = "locator"
$ws.Range("E1").Value = "locatorType"
$ws.Range("F1").Value = "value"
$ws.Range("G1").Value = "waitBefore"
$ws.Range("H1").Value = "waitAfter"
$ws.Range("A2").Value = "TC001"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "goto"
$ws.Range("D2").Value = "https://www.amazon.com/"
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 5000
$ws.Range("A3").Value = "TC001"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "waitfortext"
$ws.Range("D3").Value = "Hello, Sign in"
$ws.Range("A4").Value = "TC001"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "click"
$ws.Range("D4").Value = "Search Amazon"
$ws.Range("E4").Value = "input"
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 3000
$ws.Range("A5").Value = "TC001"
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "type"
$ws.Range("D5").Value = "Search Amazon"
$ws.Range("E5").Value = "input"
$ws.Range("F5").Value = "SAMSUNG 85-Inch Class Crystal UHD 4K DU7200 Series HDR Smart TV w/Object Tracking Sound Lite, PurColor, Motion Xcelerator, Mega Contrast, Q-Symphony (UN85DU7200, 2024 Model)"
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 3000
$ws.Range("A6").Value = "TC001"
$ws.Range("B6").Value = "Yes"
$ws.Range("C6").Value = "click"
$ws.Range("D6").Value = "Go"
$ws.Range("E6").Value = "input"
$ws.Range("G6").Value = 2000
$ws.Range("H6").Value = 4000
$ws.Range("A7").Value = "TC001"
$ws.Range("B7").Value = "Yes"
$ws.Range("C7").Value = "scroll"
$ws.Range("D7").Value = "See options"
$ws.Range("E7").Value = "a"
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 1000
$ws.Range("H7").Value = 2000
$ws.Range("A8").Value = "TC001"
$ws.Range("B8").Value = "Yes"
$ws.Range("C8").Value = "waitfortext"
$ws.Range("D8").Value = "Add to List"
$ws.Range("E8").Value = "a"
$ws.Range("G8").Value = 1000

# --- Bold style for action column data rows (reuses existing bold cellXf) ---
$ws.Range("C3:C12").Font.Bold = $true

# --- Hyperlink-like style for D13:E13 (reuse existing cellXf via copy/paste of formats) ---
$scratch = $ws.Range("Z100")
$scratch.Style = "Hyperlink"
$scratch.HorizontalAlignment = -4131
$scratch.VerticalAlignment = -4160
$scratch.Copy() | Out-Null
$ws.Range("D13:E13").PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# --- Selection ---
$ws.Range("A6").Select() | Out-Null
